$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells: C1, D1 ---
$ws.Range("C1").Value = "Retirado Por"
$ws.Range("D1").Value = "Data Retirada"

# Copy the header style (bold, bordered, centered) from B1 onto the new headers
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rows 2 and 3: voucher already used, no withdrawal info yet ---
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""

# --- Rows 4-7: vouchers withdrawn by Lucas Almeida ---
$nome = "Lucas Almeida (CPF: 28266996873)"

$ws.Range("C4").Value = $nome
$ws.Range("D4").Value = "2024-01-24 15:56:03"

$ws.Range("C5").Value = $nome
$ws.Range("D5").Value = "2024-01-24 15:56:03"

$ws.Range("C6").Value = $nome
$ws.Range("D6").Value = "2024-01-24 16:22:03"

$ws.Range("C7").Value = $nome
$ws.Range("D7").Value = "2024-01-24 17:12:57"

# --- Rows 8-18: these vouchers are now available again (status update) ---
$rows8to18 = 8..18
foreach ($r in $rows8to18) {
    $ws.Cells.Item($r, 2).Value = "Disponível"
    $ws.Cells.Item($r, 3).Value = ""
    $ws.Cells.Item($r, 4).Value = ""
}

# --- Rows 19-22: already available, just add the empty new columns ---
$rows19to22 = 19..22
foreach ($r in $rows19to22) {
    $ws.Cells.Item($r, 3).Value = ""
    $ws.Cells.Item($r, 4).Value = ""
}
